{"js": "// Translate the English interview-guide prompts to Afrikaans.\n// Each edit below locates the exact (still-English) run text with\n// context.document.body.search(...) and rewrites just that run via\n// Range.insertText(..., \"Replace\") so sibling runs (and their formatting)\n// are left untouched, mirroring the source diff which only ever swaps the\n// text inside existing <w:t> runs.\n\nasync function replaceOnce(body, findText, newText, opts) {\n  const results = body.search(findText, Object.assign({ matchCase: true }, opts));\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for \" + JSON.stringify(findText) + \" but found \" + results.items.length\n    );\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// Paragraph: \"What didn't you like about this module? How can we make it better?\"\nawait replaceOnce(\n  context.document.body,\n  \"What didn\\u2019t you like about this module? How can we make \",\n  \"Wat het jy nie van hierdie module gehou nie? Hoe kan ons \"\n);\nawait replaceOnce(\n  context.document.body,\n  \"it\",\n  \"dit\",\n  { matchWholeWord: true }\n);\nawait replaceOnce(\n  context.document.body,\n  \" better?\",\n  \" beter maak?\"\n);\n\n// Paragraph: \"What did you think of the comics, tips and examples used in the module?\"\nawait replaceOnce(\n  context.document.body,\n  \"What did you think of the \",\n  \"Wat het jy gedink van die \"\n);\nawait replaceOnce(\n  context.document.body,\n  \"comics, tips\",\n  \"strokiesprente, wenke\"\n);\nawait replaceOnce(\n  context.document.body,\n  \" and examples used in the module? \",\n  \" en voorbeelde wat in hierdie module gebruik is? \"\n);\n\n// Paragraph: \"Probe \u2013 was there anything new that this module made you think about? ...\"\nawait replaceOnce(\n  context.document.body,\n  \"Probe \\u2013 was there anything new that this module made you think about? Was daar enige nuwe aksies wat jy met jou tiener geneem het as gevolg daarvan (bv. re\\u00ebls of gesprekke, ens.)?\",\n  \"Ondersoekvraag \\u2013 was daar iets nuuts waaraan hierdie module jou laat dink het? Was daar enige nuwe aksies wat jy met jou tiener geneem het as gevolg daarvan (bv. re\\u00ebls of gesprekke, ens.)?\"\n);\n\n// Paragraph: \"Were there things which came up after doing this module which you felt more prepared to manage? Deel asseblief?\"\nawait replaceOnce(\n  context.document.body,\n  \"Were there things which came up after doing this module which you felt more prepared to manage? Deel asseblief?\",\n  \"Was daar dinge wat n\\u00e1 die voltooiing van hierdie module na vore gekom het wat jy meer voorbereid gevoel het om te hanteer? Deel asseblief?\"\n);\n", "ps1": "# Translate the English interview-guide prompts to Afrikaans.\n# Replace-DocText drives Find/Replace on $d.Content (the whole document\n# story) so only the matched run text changes - sibling runs keep their\n# own formatting/content untouched, mirroring the source diff which only\n# ever swaps the text inside existing runs.\n#\n# NOTE: all three arguments are always passed positionally - this\n# interpreter does not reliably apply boolean parameter defaults.\n\n$d = $word.ActiveDocument\n\nfunction Replace-DocText($FindText, $ReplaceText, $WholeWord) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $FindText\n    $find.Replacement.Text = $ReplaceText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $WholeWord\n    $find.MatchWildcards = $false\n\n    $ok = $find.Execute(\n        $find.Text,\n        $true,\n        $WholeWord,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $find.Replacement.Text,\n        2\n    )\n\n    if (-not $ok) {\n        throw \"Find/Replace did not find a match for: $FindText\"\n    }\n}\n\n# Paragraph: \"What didn\u2019t you like about this module? How can we make it better?\"\nReplace-DocText \"What didn\u2019t you like about this module? How can we make \" \"Wat het jy nie van hierdie module gehou nie? Hoe kan ons \" $false\nReplace-DocText \"it\" \"dit\" $true\nReplace-DocText \" better?\" \" beter maak?\" $false\n\n# Paragraph: \"What did you think of the comics, tips and examples used in the module?\"\nReplace-DocText \"What did you think of the \" \"Wat het jy gedink van die \" $false\nReplace-DocText \"comics, tips\" \"strokiesprente, wenke\" $false\nReplace-DocText \" and examples used in the module? \" \" en voorbeelde wat in hierdie module gebruik is? \" $false\n\n# Paragraph: \"Probe \u2013 was there anything new that this module made you think about? ...\"\nReplace-DocText \"Probe \u2013 was there anything new that this module made you think about? Was daar enige nuwe aksies wat jy met jou tiener geneem het as gevolg daarvan (bv. re\u00ebls of gesprekke, ens.)?\" \"Ondersoekvraag \u2013 was daar iets nuuts waaraan hierdie module jou laat dink het? Was daar enige nuwe aksies wat jy met jou tiener geneem het as gevolg daarvan (bv. re\u00ebls of gesprekke, ens.)?\" $false\n\n# Paragraph: \"Were there things which came up after doing this module which you felt more prepared to manage? Deel asseblief?\"\nReplace-DocText \"Were there things which came up after doing this module which you felt more prepared to manage? Deel asseblief?\" \"Was daar dinge wat n\u00e1 die voltooiing van hierdie module na vore gekom het wat jy meer voorbereid gevoel het om te hanteer? Deel asseblief?\" $false\n\nWrite-Output \"done\"\n"}
